$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated simulation results for the 380 kV case (columns B,C,D,E,G,J,O; rows 2-25)
$data = @(
    @(1.499813420424744, 0.3201417225299679, 0.6567851668096125, 0.2679783297317755, 0.00247545144158739, 0.1395546008024766, 4.353606491903719),
    @(1.366108778815885, 0.2836510130846364, 0.6466607056188138, 0.2629821746456997, 0.002479536462328296, 0.1362230728056133, 4.371756388375871),
    @(1.284229195812884, 0.2612170243696426, 0.6407968441174603, 0.260065007182142, 0.002482176148803917, 0.1342605386821774, 4.386455633701189),
    @(1.250917880320003, 0.2520681751275902, 0.6384958322543923, 0.2589140082089543, 0.002483285009263252, 0.1334816071572646, 4.393336543625423),
    @(1.24538993258642, 0.2505486148288298, 0.6381190967868804, 0.2587251648846021, 0.002483471141074313, 0.1333535216403661, 4.394532827327453),
    @(1.283779722506438, 0.2610936669003365, 0.6407654534277469, 0.260049331569391, 0.002482190968918457, 0.1342499495339595, 4.386544829215808),
    @(1.453667895635874, 0.3075658530917735, 0.6532210020814944, 0.2662243748787176, 0.002476832735769057, 0.1383886157251339, 4.359124820092035),
    @(1.788500980299091, 0.3984602469055289, 0.6804509688379881, 0.2795319773607261, 0.002467363437908815, 0.1471668792328913, 4.333707972221987),
    @(2.035519218576781, 0.4650882445798743, 0.7021795208968626, 0.2900470319503086, 0.002461032294642033, 0.1540257883277576, 4.332527206724876),
    @(2.148115209504283, 0.4953650774198763, 0.7124415091518017, 0.2949925579918045, 0.002458286528881878, 0.1572363134128381, 4.335833861061758),
    @(2.19078439163718, 0.5068252250595719, 0.7163819391728623, 0.2968887330727981, 0.002457265978574328, 0.1584651407708719, 4.337642339193508),
    @(2.181593433243734, 0.5043573064216389, 0.7155308735763128, 0.2964793151284439, 0.002457484919640294, 0.1581999082505803, 4.337228051279055),
    @(2.151625002585774, 0.496308013459668, 0.7127645988415452, 0.2951480877222963, 0.00245820218318812, 0.1573371473748324, 4.335971475560143),
    @(2.133272549906735, 0.4913769236894154, 0.7110772694744014, 0.2943357241313862, 0.002458644026859331, 0.1568103859491004, 4.33527434267404),
    @(2.028165270515274, 0.4631088968653785, 0.7015164856182707, 0.2897271013134386, 0.002461214430879285, 0.1538177979376485, 4.332388762272757),
    @(1.963742561583786, 0.4457587835137815, 0.6957480586275437, 0.2869414616099419, 0.002462825616584755, 0.1520051422446045, 4.331605586459943),
    @(1.926709667324587, 0.4357764109032018, 0.6924657408931694, 0.2853544944381454, 0.002463764975308674, 0.1509710523624506, 4.331516704858956),
    @(1.914174662491405, 0.4323960386215049, 0.6913605022788829, 0.2848197926429279, 0.002464085201278499, 0.1506223851914399, 4.331548603770983),
    @(1.970598267239836, 0.4476060505206192, 0.6963584393859321, 0.2872364177976436, 0.002462652795014251, 0.1521972220956087, 4.331651508833573),
    @(2.160426609318904, 0.498672425196105, 0.7135756421839119, 0.2955384653579571, 0.002457990984991166, 0.1575902058001333, 4.336325435270453),
    @(2.284673714840267, 0.5320177703034688, 0.7251454549025027, 0.3011008383964295, 0.002455056156001162, 0.1611910573036965, 4.34262445150415),
    @(2.218344242929447, 0.5142235415702316, 0.7189413392204642, 0.2981195743667371, 0.002456612320713801, 0.1592622151727454, 4.338964501427142),
    @(1.967498788374996, 0.4467709246371783, 0.6960823802671712, 0.2871030228416842, 0.002462730886951106, 0.1521103577895389, 4.331629621940493),
    @(1.697740725735741, 0.3738971273527341, 0.6727829641563972, 0.2758028668630246, 0.00246981471047742, 0.1447206289294769, 4.337527014933158)
)

$cols = @("B", "C", "D", "E", "G", "J", "O")
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $colLetter = $cols[$j]
        $cellRef = $colLetter + $rowNum
        $ws.Range($cellRef).Value = $data[$i][$j]
    }
}
